$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 20374.56648578441,  95.01934733259309,  0.004641977033988871, 2019),
    @(2, 43640.44824045801,  1030.429758063137,   0.02306714808912531,  2019),
    @(3, 8871.998226522564,  1916.056733544713,   0.1776091001238989,   2019),
    @(4, 1562.395414645924,  1945.754422583905,   0.5546383458126031,   2019)
)

$startRow = 34
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}
